$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet's columns are being re-ordered (and four new fields are being
# added: Initials, Title, Country, EmailAddress). Rather than shuffle
# existing columns around, just rewrite row 1 (headers) and row 2 (data)
# directly in their final target layout, column by column.
# ---------------------------------------------------------------------------

# --- Row 1: headers -----------------------------------------------------
$ws.Range("A1").Value = "First"
$ws.Range("B1").Value = "Last"
$ws.Range("C1").Value = "Username"
$ws.Range("D1").Value = "Initials"
$ws.Range("E1").Value = "Title"
$ws.Range("F1").Value = "Department"
$ws.Range("G1").Value = "Company"
$ws.Range("H1").Value = "Office"
$ws.Range("I1").Value = "StreetAddress"
$ws.Range("J1").Value = "POBox"
$ws.Range("K1").Value = "City"
$ws.Range("L1").Value = "State"
$ws.Range("M1").Value = "Country"
$ws.Range("N1").Value = "PostalCode"
$ws.Range("O1").Value = "Organization"
$ws.Range("P1").Value = "EmployeeNumber"
$ws.Range("Q1").Value = "EmployeeID"
$ws.Range("R1").Value = "OfficePhone"
$ws.Range("S1").Value = "MobilePhone"
$ws.Range("T1").Value = "HomePhone"
$ws.Range("U1").Value = "Fax"
$ws.Range("V1").Value = "Manager"
$ws.Range("W1").Value = "EmailAddress"

# --- Row 2: data ----------------------------------------------------------
$ws.Range("A2").Value = "John"
$ws.Range("B2").Value = "Smith"
$ws.Range("C2").Formula = "=LEFT(A2,1)&B2"
$ws.Range("D2").Formula = "=LEFT(A2,1)&LEFT(B2,1)"
$ws.Range("E2").Value = "Self-Supervisor"
$ws.Range("F2").Value = "Operations"
$ws.Range("G2").Value = "Company 1"
$ws.Range("H2").Value = 101
$ws.Range("I2").Value = "123 Fake ST"
$ws.Range("J2").Value = 101
$ws.Range("K2").Value = "San Jose"
$ws.Range("L2").Value = "CA"
$ws.Range("M2").Value = "USA"
$ws.Range("N2").Value = 8675309
$ws.Range("O2").Value = "Unit 1"
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = "123-456-7892"
$ws.Range("S2").Value = "123-456-789"
$ws.Range("T2").Value = "123-456-7893"
$ws.Range("U2").Value = "123-456-7891"
$ws.Range("V2").Value = "jsmith"
$ws.Range("W2").Value = "john.smith@company.com"

# EmailAddress cell becomes a live mailto: hyperlink (Excel auto-applies the
# built-in "Hyperlink" cell style -- underlined, theme color 10).
$ws.Hyperlinks.Add($ws.Range("W2"), "mailto:john.smith@company.com")

# Auto-size the new/changed columns to fit their content, matching Excel's
# behaviour when columns are populated interactively.
$ws.Columns("D:D").AutoFit()
$ws.Columns("F:F").AutoFit()
$ws.Columns("W:W").AutoFit()

# Restore the active cell/selection seen in the saved workbook.
$ws.Range("E11").Select() | Out-Null
